$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates (per diff) ---
$ws.Range("D2").Value = "28.378.87"
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").Value = "1.805.71"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "316.06"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.5499"
$ws.Range("E7").Value = "  +6.03%  "
$ws.Range("D8").Value = "0.3853"
$ws.Range("E8").Value = "  +6.64%  "
$ws.Range("D9").Value = "0.07593"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").Value = "42.56"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +2.87%  "
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "21.17"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "6.196"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "7.379"
$ws.Range("E15").Value = "  +5.29%  "
$ws.Range("D16").Value = "1.800.72"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "92.25"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "0.06443"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "0.9988"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").Value = "5.989"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "28.397.46"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "2.131"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").Value = "158.32"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").Value = "20.65"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "2.403"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "2.012.68"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "123.70"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "1.124"
$ws.Range("E31").Value = "  +5.47%  "
$ws.Range("D32").Value = "0.1022"
$ws.Range("E32").Value = "  +5.18%  "
$ws.Range("D33").Value = "5.750"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("D35").Value = "0.2337"
$ws.Range("E35").Value = "  +15.03%  "
$ws.Range("D36").Value = "0.06381"
$ws.Range("E36").Value = "  +6.35%  "
$ws.Range("D37").Value = "0.02325"
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("D38").Value = "8.853"
$ws.Range("E38").Value = "  +10.10%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.109"
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "11.64"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").Value = "0.6426"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("D42").Value = "1.162"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").Value = "0.9986"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "1.382"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "13.50"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "0.5979"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").Value = "126.22"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").Value = "0.06906"
$ws.Range("E51").Value = "  +2.74%  "
